$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix player-name typos (trailing spaces / abbreviated prefixes) found via
# the "find broken player names" pass referenced in the commit message.
$ws.Range("S6").Value = "PELLEGRINI"
$ws.Range("K7").Value = "LOPEZ"
$ws.Range("C10").Value = "FARES"
$ws.Range("B4").Value = "MARRONE"

# Restore the previously-selected cell
[void]$ws.Range("B4").Select()
